$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "248.05"
Set-TextValue "E2" "1BNBBNBBestin24h"

# Row 3 - OKB
Set-TextValue "D3" "22.40"

# Row 5 - Cronos
Set-TextValue "D5" "0.05611"

# Row 6 - GateToken
Set-TextValue "D6" "3.400"

# Row 7 - KuCoinToken
Set-TextValue "D7" "6.467"

# Row 8 - was FTXToken, now MXToken
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.8022"
Set-TextValue "E8" "7MXTokenMX"

# Row 9 - was MXToken, now FTXToken
Set-TextValue "B9" "FTXToken"
Set-TextValue "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "1.066"
Set-TextValue "E9" "8FTXTokenFTT"

# Row 10 - WazirX
Set-TextValue "D10" "0.1426"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextValue "D12" "0.03192"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.02967"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09266"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001670"

# Row 16 - MCDex
Set-TextValue "D16" "3.251"

# Row 17 - CoinExToken
Set-TextValue "D17" "0.04695"

# Row 18 - One
Set-TextValue "D18" "0.0005742"
Set-TextValue "E18" "17OneONE"

# Row 19 - TigerCash
Set-TextValue "D19" "0.006268"

# Row 20 - BitKan
Set-TextValue "D20" "0.001048"

# Row 22 - NitroEx
Set-TextValue "D22" "0.0001501"

# Row 24 - LEO
Set-TextValue "D24" "3.980"

# Row 27 - ProBitToken
Set-TextValue "D27" "0.1277"

# Row 40 - IDEX
Set-TextValue "D40" "0.04191"

# Row 41 - KickToken
Set-TextValue "D41" "0.003250"
Set-TextValue "E41" "40KickTokenKICKWorstin24h"

# Row 42 - was CEJI, now BKEXToken
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1047"
Set-TextValue "E42" "41BKEXTokenBKK"

# Row 43 - was BKEXToken, now CEJI
Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002971"
Set-TextValue "E43" "42CEJICEJI"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008691"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005660"

# Row 48 - BOLO
Set-TextValue "D48" "0.02930"
Set-TextValue "E48" "47BOLOBOLO"

# Row 49 - CryptobidCoin
Set-TextValue "D49" "0.00002101"
